# "Fruta / hortaliza, semanal" weekly update.
#
# This adds a new week of price reports (5 new rows) for the
# "Vega Modelo de Temuco - Pimiento" subset. The new rows are inserted
# right before the current row 1162, pushing the existing rows
# 1162-1175 down to 1167-1180 (dimension grows from R1175 to R1180).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows above row 1162 (existing data shifts down).
$ws.Range("A1162:R1166").Insert()

# Columns A, B, C, E, F, G, R are constant for every record in this
# subset (market/region/category metadata).
$commonCols = @{
    A = 10
    B = 'Vega Modelo de Temuco'
    C = 'La Araucanía'
    E = 9
    F = 100112002
    G = 'Pimiento'
    R = 'Hortaliza'
}

# New weekly records for rows 1162-1166.
$newRecords = @(
    @{ Row=1162; D=44595; H='Cuatro cascos amarillo'; I='Extra';    J=40;  K=25000; L=25000; M=25000; N='$/caja 15 kilos'; O='Región del Maule';             P=1667; Q=15 }
    @{ Row=1163; D=44595; H='Cuatro cascos amarillo'; I='Primera';  J=100; K=20000; L=20000; M=20000; N='$/caja 15 kilos'; O='Región del Maule';             P=1333; Q=15 }
    @{ Row=1164; D=44595; H='Cuatro cascos verde';    I='Primera';  J=400; K=10000; L=10000; M=10000; N='$/caja 15 kilos'; O='Región del Maule';             P=667;  Q=15 }
    @{ Row=1165; D=44595; H='Zafiro rojo';            I='Primera';  J=300; K=20000; L=20000; M=20000; N='$/caja 15 kilos'; O='Región de Arica y Parinacota'; P=1333; Q=15 }
    @{ Row=1166; D=44595; H='Zafiro rojo';            I='Segunda';  J=80;  K=15000; L=15000; M=15000; N='$/caja 15 kilos'; O='Región de Arica y Parinacota'; P=1000; Q=15 }
)

foreach ($rec in $newRecords) {
    $r = $rec.Row

    foreach ($col in $commonCols.Keys) {
        $ws.Range("$col$r").Value = $commonCols[$col]
    }

    $ws.Range("D$r").Value = $rec.D
    $ws.Range("H$r").Value = $rec.H
    $ws.Range("I$r").Value = $rec.I
    $ws.Range("J$r").Value = $rec.J
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("N$r").Value = $rec.N
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
    $ws.Range("Q$r").Value = $rec.Q
}
